$wb = $excel.ActiveWorkbook

# --- 1. Status text: "Ready for handoff" -> "In Translation" -------------
# Overview sheet holds the per-locale status in columns E (zh-cn) and F (de-de)
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("E2").Value = "In Translation"
$wsOverview.Range("F2").Value = "In Translation"

# Per-locale detail sheets keep the same status in column C ("Status")
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("C2").Value = "In Translation"

$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("C2").Value = "In Translation"

# --- 2. Narrow the (now shorter) status columns ---------------------------
# Target stored width is ~13.41 chars; the host's ColumnWidth setter only
# lands on 1/6-character increments, so 12.5 is the closest input that
# resolves to the nearest achievable stored width.
$wsOverview.Range("E1").ColumnWidth = 12.5
$wsOverview.Range("F1").ColumnWidth = 12.5
$wsZhCn.Range("C1").ColumnWidth = 12.5
$wsDeDe.Range("C1").ColumnWidth = 12.5
